# Update "想去人数" (want-to-go count) values in column F across all sheets.
# This mirrors the upstream data refresh captured in commit:
#   "Update gh-pages to output generated at 456a3b4"
#
# Sheets affected (by name):
#   展览      (Exhibitions)
#   演出      (Performances)
#   本地生活  (Local life)
#   全部类型  (All types - combined view)

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 11505  # F2: 11502 -> 11505
$ws.Cells.Item(4, 6).Value = 627  # F4: 626 -> 627
$ws.Cells.Item(5, 6).Value = 485  # F5: 483 -> 485
$ws.Cells.Item(6, 6).Value = 1442  # F6: 1434 -> 1442
$ws.Cells.Item(7, 6).Value = 734  # F7: 728 -> 734
$ws.Cells.Item(8, 6).Value = 165  # F8: 160 -> 165
$ws.Cells.Item(9, 6).Value = 36  # F9: 33 -> 36
$ws.Cells.Item(10, 6).Value = 1044  # F10: 1041 -> 1044
$ws.Cells.Item(11, 6).Value = 613  # F11: 605 -> 613
$ws.Cells.Item(12, 6).Value = 718  # F12: 714 -> 718
$ws.Cells.Item(13, 6).Value = 1231  # F13: 1217 -> 1231
$ws.Cells.Item(14, 6).Value = 247  # F14: 241 -> 247
$ws.Cells.Item(15, 6).Value = 975  # F15: 973 -> 975
$ws.Cells.Item(16, 6).Value = 35  # F16: 34 -> 35
$ws.Cells.Item(17, 6).Value = 168  # F17: 162 -> 168
$ws.Cells.Item(18, 6).Value = 22  # F18: 21 -> 22
$ws.Cells.Item(19, 6).Value = 356  # F19: 352 -> 356
$ws.Cells.Item(21, 6).Value = 277  # F21: 274 -> 277
$ws.Cells.Item(22, 6).Value = 502  # F22: 499 -> 502
$ws.Cells.Item(23, 6).Value = 526  # F23: 523 -> 526
$ws.Cells.Item(24, 6).Value = 716  # F24: 712 -> 716
$ws.Cells.Item(26, 6).Value = 136  # F26: 135 -> 136

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 939  # F4: 936 -> 939
$ws.Cells.Item(6, 6).Value = 156  # F6: 155 -> 156
$ws.Cells.Item(7, 6).Value = 11  # F7: 10 -> 11
$ws.Cells.Item(10, 6).Value = 482  # F10: 471 -> 482

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 145  # F2: 141 -> 145

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 11505  # F2: 11502 -> 11505
$ws.Cells.Item(5, 6).Value = 627  # F5: 626 -> 627
$ws.Cells.Item(6, 6).Value = 485  # F6: 483 -> 485
$ws.Cells.Item(7, 6).Value = 145  # F7: 141 -> 145
$ws.Cells.Item(8, 6).Value = 1442  # F8: 1434 -> 1442
$ws.Cells.Item(10, 6).Value = 734  # F10: 728 -> 734
$ws.Cells.Item(11, 6).Value = 165  # F11: 160 -> 165
$ws.Cells.Item(12, 6).Value = 939  # F12: 936 -> 939
$ws.Cells.Item(13, 6).Value = 36  # F13: 33 -> 36
$ws.Cells.Item(14, 6).Value = 1044  # F14: 1041 -> 1044
$ws.Cells.Item(15, 6).Value = 613  # F15: 605 -> 613
$ws.Cells.Item(16, 6).Value = 718  # F16: 714 -> 718
$ws.Cells.Item(17, 6).Value = 1231  # F17: 1217 -> 1231
$ws.Cells.Item(18, 6).Value = 247  # F18: 241 -> 247
$ws.Cells.Item(19, 6).Value = 975  # F19: 973 -> 975
$ws.Cells.Item(20, 6).Value = 35  # F20: 34 -> 35
$ws.Cells.Item(21, 6).Value = 168  # F21: 162 -> 168
$ws.Cells.Item(22, 6).Value = 22  # F22: 21 -> 22
$ws.Cells.Item(23, 6).Value = 356  # F23: 352 -> 356
$ws.Cells.Item(25, 6).Value = 156  # F25: 155 -> 156
$ws.Cells.Item(27, 6).Value = 277  # F27: 274 -> 277
$ws.Cells.Item(28, 6).Value = 11  # F28: 10 -> 11
$ws.Cells.Item(31, 6).Value = 502  # F31: 499 -> 502
$ws.Cells.Item(32, 6).Value = 526  # F32: 523 -> 526
$ws.Cells.Item(33, 6).Value = 716  # F33: 712 -> 716
$ws.Cells.Item(36, 6).Value = 136  # F36: 135 -> 136
$ws.Cells.Item(37, 6).Value = 482  # F37: 471 -> 482

